# Actualización desde MV -datos-
# Adds the daily "Tasas de Bonos de Gobierno a 10 años" data for
# 07-09-2021 .. 13-09-2021 and corrects the previously-provisional
# 06-09-2021 (row 178) figures now that final data is available.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($row, $colLetter, $text) {
    $cell = $ws.Cells.Item($row, (Col-Num $colLetter))
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Col-Num($col) {
    $n = 0
    foreach ($ch in $col.ToCharArray()) {
        $n = $n * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $n
}

function Set-Row($rowNum, $dateText, $values) {
    Set-DateText $rowNum "A" $dateText
    foreach ($col in $values.Keys) {
        $ws.Cells.Item($rowNum, (Col-Num $col)).Value = $values[$col]
    }
}

# --- Correct the existing row 178 (06-09-2021) now final data is in ---
$ws.Cells.Item(178, (Col-Num "C")).Value = 0.69
$ws.Cells.Item(178, (Col-Num "E")).Value = -0.37
$ws.Cells.Item(178, (Col-Num "H")).Value = 1.93
$ws.Cells.Item(178, (Col-Num "I")).Value = 1.78
$ws.Cells.Item(178, (Col-Num "L")).Value = 6.04
$ws.Cells.Item(178, (Col-Num "M")).Value = 10.83
$ws.Cells.Item(178, (Col-Num "N")).Value = 7.28
$ws.Cells.Item(178, (Col-Num "O")).Value = 6.98
$ws.Cells.Item(178, (Col-Num "P")).Value = 6.29

# --- New rows 179-183 ---
Set-Row 179 "07-09-2021" @{
    B = 1.37; C = 0.74; D = 0.04; E = -0.32; F = 1.99; G = 3.22
    H = 1.99; I = 1.79; J = 7.03; K = 1.66; L = 6.06
    N = 7.32; O = 7.02; P = 6.29
}

Set-Row 180 "08-09-2021" @{
    B = 1.34; C = 0.74; D = 0.05; E = -0.32; F = 2.03; G = 3.23
    H = 2; I = 1.81; J = 7.03; K = 1.71; L = 6.09
    M = 11.1; N = 7.31; O = 6.98; P = 6.38
}

Set-Row 181 "09-09-2021" @{
    B = 1.3; C = 0.74; D = 0.04; E = -0.36; F = 2.03; G = 3.22
    H = 1.95; I = 1.84; J = 7.02; K = 1.7; L = 6.08
    M = 10.87; N = 7.38; O = 6.94; P = 6.41
}

Set-Row 182 "10-09-2021" @{
    B = 1.34; C = 0.76; D = 0.05; E = -0.33; F = 2.02; G = 3.25
    H = 1.99; I = 1.88; J = 7.01; K = 1.67; L = 6.08
    M = 11.04; N = 7.39; O = 6.95; P = 6.44
}

Set-Row 183 "13-09-2021" @{
    B = 1.33; C = 0.74; D = 0.05; E = -0.34; F = 2.06; G = 3.26
    H = 1.99; I = 1.89; J = 6.99; K = 1.68; L = 6.11
    M = 11.01; O = 6.95
}
